$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10801997.17295967
$ws.Range("C2").Value = 2374074.297282989
$ws.Range("D2").Value = 26378827.89743374
$ws.Range("E2").Value = 1185682.38352715
$ws.Range("F2").Value = 8439986.56294816
$ws.Range("G2").Value = 1828197.521994759
$ws.Range("H2").Value = 2104549.666711418
$ws.Range("I2").Value = 10801997.17295967
$ws.Range("J2").Value = 45512584
$ws.Range("L2").Value = 28752902.19471672
$ws.Range("M2").Value = 9625668.94647531
$ws.Range("N2").Value = 3932747.188706177
$ws.Range("O2").Value = 45117.65828076528
$ws.Range("P2").Value = 206824.7031532861
$ws.Range("Q2").Value = 251942.3614340514
$ws.Range("S2").Value = 116051.4494870003
$ws.Range("T2").Value = 116051.4494870003
